# ---------------------------------------------------------------------------
# Adds the "kPROT" hydrophobicity-scale columns (P:S) to Tabelle1, matching
# the upstream commit that introduced the Pilpel/Ben-Tal/Lancet KPROT scale
# (kPROT_Extracellular, kPROT_Central, kPROT_Intracellular,
# kPROT_Both termini) alongside the existing scales in columns A:O.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Row 2 (merged-looking "Reference" row, id column G/N/O etc.) - the four new
# columns all cite the same KPROT reference, exactly like N2/O2 already do.
$reference = "Pilpel Y, Ben-Tal N, & Lancet D (1999) KPROT: A knowledge-based scale for the propensity of residue orientation in transmembrane segments. Application to membrane protein structure prediction11Edited by G. von Heijne. J. Mol. Biol. 294(4):921-935."
$ws.Range("P2").Value = $reference
$ws.Range("Q2").Value = $reference
$ws.Range("R2").Value = $reference
$ws.Range("S2").Value = $reference

# Row 3: column headers for the four new KPROT scales. Copy O3's formatting
# (bold/centered header style) onto the new header cells first.
$ws.Range("O3").Copy() | Out-Null
$ws.Range("P3:S3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("P3").Value = "kPROT_Extracellular"
$ws.Range("Q3").Value = "kPROT_Central"
$ws.Range("R3").Value = "kPROT_Intracellular"
$ws.Range("S3").Value = "kPROT_Both termini"

# Rows 4-23: the KPROT values per amino acid, for the four new scale columns.
$kprotData = @{
    4  = @(0.03, 0.09, 0.18, 0.09)
    5  = @(0.53, 0.12, 0.61, 0.5)
    6  = @(1.35, 1.1,  1.23, 0.73)
    7  = @(1.01, 1.1,  1.33, 0.7)
    8  = @(0.41, 0.16, 0.12, 0.07)
    9  = @(0.02, 0.05, 0.33, 0.23)
    10 = @(0.21, 0.69, 0.21, 0.15)
    11 = @(0.34, 0.12, 0.09, 0.13)
    12 = @(0.85, 0.56, 0.66, 0.46)
    13 = @(0.11, 0.26, 0.17, 0.09)
    14 = @(0.28, 0.39, 0.32, 0.22)
    15 = @(0.69, 1.11, 0.75, 0.38)
    16 = @(0.16, 0.66, 0.75, 0.37)
    17 = @(0.78, 0.83, 0.39, 0.33)
    18 = @(0.53, 0.84, 0.44, 0.4)
    19 = @(0.11, 0.22, 0.41, 0.21)
    20 = @(0.03, 0.03, 0.43, 0)
    21 = @(0.27, 0.31, 0,    0.11)
    22 = @(0.25, 0.65, 0.29, 0.16)
    23 = @(0.18, 0.7,  0.26, 0.23)
}

foreach ($r in 4..23) {
    $vals = $kprotData[$r]
    $ws.Cells.Item($r, 16).Value = $vals[0]   # P
    $ws.Cells.Item($r, 17).Value = $vals[1]   # Q
    $ws.Cells.Item($r, 18).Value = $vals[2]   # R
    $ws.Cells.Item($r, 19).Value = $vals[3]   # S
}

# Restore the selection to match the post-edit workbook (user ended up with
# column P in view, P15 selected).
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("P15").Select() | Out-Null
